$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    2.2902340000000501,
    4.5494420000031797,
    6.8264019999915,
    9.0982269999949494,
    11.3923329999961,
    13.701493000000401,
    16.0234329999948,
    18.3126649999903,
    20.651316999996101,
    22.944239999996999,
    25.2865619999938,
    27.626585999998401,
    29.9653560000006,
    32.409829999989597,
    34.891206999993301,
    37.461176999990101,
    39.962643999999202,
    42.508680999992002,
    45.137189999994,
    47.831139000001698,
    50.588176999997799,
    53.373261999993701,
    56.205810999992501,
    59.113016000002901,
    62.081175999992404,
    65.053748999998703,
    68.105171999995903,
    71.249364999996004,
    74.452296999996094,
    77.748700999989495
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

$ws.Range("G18").Select()
